$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data set (id, building, roomNo, capacity) for rows 2..26
$data = @(
    @(36, "Old Acad", "C24",  18),
    @(37, "LHC",      "C101", 169),
    @(38, "LHC",      "C102", 281),
    @(39, "LHC",      "C201", 177),
    @(40, "LHC",      "C208", 20),
    @(41, "LHC",      "C209", 22),
    @(42, "LHC",      "C210", 25),
    @(43, "Old Acad", "C211", 24),
    @(44, "LHC",      "C212", 28),
    @(45, "LHC",      "C213", 31),
    @(46, "LHC",      "C214", 20),
    @(47, "LHC",      "C215", 22),
    @(48, "LHC",      "C216", 17),
    @(49, "RnD",      "A006", 58),
    @(50, "RnD",      "A007", 56),
    @(51, "RnD",      "A106", 48),
    @(52, "RnD",      "B105", 42),
    @(53, "Old Acad", "C01",  100),
    @(54, "Old Acad", "C02",  100),
    @(55, "Old Acad", "C03",  100),
    @(56, "Old Acad", "C11",  100),
    @(57, "Old Acad", "C04",  100),
    @(58, "RnD",      "A105", 100),
    @(59, "Old Acad", "C12",  100),
    @(60, "Old Acad", "C21",  100)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $entry = $data[$i]
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
}
